# Update cryptos list (simulated "Updated cryptos list ... with GitHub Actions" run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written as TEXT (avoids Excel's automatic
# number coercion for numeric-looking strings like "224.51"), while
# restoring the cell's original (unstyled) appearance afterwards.
function Set-TextValue {
    param($Worksheet, [string]$Address, [string]$Text)
    $r = $Worksheet.Range($Address)
    $r.NumberFormat = "@"
    $r.Value = $Text
    $r.Style = "Normal"
}

# --- Row 2: Bitcoin ---
Set-TextValue $ws "D2" "33.861.20"
$ws.Range("E2").Value = "  -1.05%  "

# --- Row 3: Ethereum ---
Set-TextValue $ws "D3" "1.781.43"
$ws.Range("E3").Value = "  -1.44%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.07%  "

# --- Row 5: BNB ---
Set-TextValue $ws "D5" "224.51"
$ws.Range("E5").Value = "  +0.57%  "

# --- Row 6: XRP ---
$ws.Range("E6").Value = "  -1.41%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  +0.06%  "

# --- Row 8: Solana ---
Set-TextValue $ws "D8" "31.77"
$ws.Range("E8").Value = "  -3.81%  "

# --- Row 9: Cardano ---
$ws.Range("E9").Value = "  +0.22%  "

# --- Row 10: Dogecoin ---
$ws.Range("E10").Value = "  -5.60%  "

# --- Row 11: TRON ---
Set-TextValue $ws "D11" "0.0935"
$ws.Range("E11").Value = "  +0.74%  "

# --- Row 12: WrappedliquidstakedEther2.0 ---
Set-TextValue $ws "D12" "2.037.10"
$ws.Range("E12").Value = "  -1.48%  "

# --- Row 13: Chainlink ---
Set-TextValue $ws "D13" "11.14"
$ws.Range("E13").Value = "  +0.65%  "

# --- Row 14: WrappedEther ---
Set-TextValue $ws "D14" "1.787.05"
$ws.Range("E14").Value = "  -1.28%  "

# --- Row 15: WrappedBTC ---
Set-TextValue $ws "D15" "33.860.40"
$ws.Range("E15").Value = "  -1.18%  "

# --- Row 16: Polygon ---
$ws.Range("E16").Value = "  -3.63%  "

# --- Row 17: Polkadot ---
$ws.Range("E17").Value = "  -2.42%  "

# --- Row 18: Litecoin ---
Set-TextValue $ws "D18" "66.89"
$ws.Range("E18").Value = "  -2.91%  "

# --- Row 19: BitcoinCash ---
Set-TextValue $ws "D19" "239.11"
$ws.Range("E19").Value = "  -3.52%  "

# --- Row 20: ShibaInu ---
$ws.Range("E20").Value = "  -2.23%  "

# --- Row 22: Avalanche ---
Set-TextValue $ws "D22" "10.55"
$ws.Range("E22").Value = "  -4.90%  "

# --- Row 23: Uniswap ---
$ws.Range("E23").Value = "  -2.38%  "

# --- Row 24: Toncoin ---
$ws.Range("E24").Value = "  -3.26%  "

# --- Row 25: Monero ---
Set-TextValue $ws "D25" "160.92"

# --- Rows 26 & 27: swap Cosmos/EthereumClassic with updated values ---
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D26" "16.10"
$ws.Range("E26").Value = "  -3.14%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws "D27" "7.02"
$ws.Range("E27").Value = "  -2.08%  "

# --- Row 28: Stellar ---
$ws.Range("E28").Value = "  -0.98%  "

# --- Row 29: BinanceUSD ---
$ws.Range("E29").Value = "  +0.18%  "

# --- Row 30: PancakeSwap ---
$ws.Range("E30").Value = "  +0.83%  "

# --- Row 32: Filecoin ---
$ws.Range("E32").Value = "  -4.12%  "

# --- Row 33: InternetComputer(DFINITY) ---
$ws.Range("E33").Value = "  -0.51%  "

# --- Row 34: LidoDAOToken ---
$ws.Range("E34").Value = "  -2.35%  "

# --- Row 35: Maker ---
Set-TextValue $ws "D35" "1.390.34"
$ws.Range("E35").Value = "  -2.18%  "

# --- Row 36: ImmutableX ---
Set-TextValue $ws "D36" "0.636"
$ws.Range("E36").Value = "  -2.67%  "

# --- Row 37: TrustWalletToken ---
$ws.Range("E37").Value = "  -1.79%  "

# --- Row 38: VeChain ---
$ws.Range("E38").Value = "  -1.28%  "

# --- Row 39: RenderToken ---
Set-TextValue $ws "D39" "2.26"
$ws.Range("E39").Value = "  +3.98%  "

# --- Row 40: HuobiToken ---
$ws.Range("E40").Value = "  -0.02%  "

# --- Row 41: ARBITRUM ---
$ws.Range("E41").Value = "  -3.46%  "

# --- Row 42: Aave ---
Set-TextValue $ws "D42" "78.32"
$ws.Range("E42").Value = "  -3.23%  "

# --- Rows 43 & 44: swap MXToken/InjectiveProtocol with updated values ---
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws "D43" "13.52"
$ws.Range("E43").Value = "  +11.90%  "

$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D44" "2.65"
$ws.Range("E44").Value = "  -3.13%  "

# --- Row 45: Kaspa ---
$ws.Range("E45").Value = "  +2.49%  "

# --- Row 46: WEMIXToken ---
$ws.Range("E46").Value = "  +2.08%  "

# --- Row 47: BabyDogeCoin ---
$ws.Range("E47").Value = "  +7.69%  "

# --- Row 48: FraxShare ---
$ws.Range("E48").Value = "  -1.79%  "

# --- Row 49: RocketPoolETH ---
Set-TextValue $ws "D49" "1.938.09"
$ws.Range("E49").Value = "  -1.55%  "

# --- Row 50: Quant ---
Set-TextValue $ws "D50" "105.15"
$ws.Range("E50").Value = "  -2.96%  "

# --- Row 51: PaxDollar ---
$ws.Range("E51").Value = "  +0.12%  "
